$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 760
$ws.Range("F4").Value = 2049
$ws.Range("F5").Value = 6104
$ws.Range("F6").Value = 3521
$ws.Range("F7").Value = 703
$ws.Range("F8").Value = 61
$ws.Range("F9").Value = 1464
$ws.Range("F10").Value = 4876
$ws.Range("F11").Value = 1123
$ws.Range("F12").Value = 1818
$ws.Range("F14").Value = 71
$ws.Range("F15").Value = 79
$ws.Range("F16").Value = 227
$ws.Range("F19").Value = 349
$ws.Range("F22").Value = 102
$ws.Range("F23").Value = 12
$ws.Range("F24").Value = 229
$ws.Range("F25").Value = 118
$ws.Range("F26").Value = 47
$ws.Range("F27").Value = 1172
$ws.Range("F28").Value = 451
$ws.Range("F29").Value = 137
$ws.Range("F30").Value = 255
$ws.Range("F31").Value = 512
$ws.Range("F32").Value = 1009
$ws.Range("F33").Value = 38
$ws.Range("F34").Value = 1875
$ws.Range("F35").Value = 2335
$ws.Range("F36").Value = 1107
$ws.Range("F38").Value = 47
$ws.Range("F39").Value = 309
$ws.Range("F40").Value = 175
$ws.Range("F41").Value = 688
$ws.Range("F42").Value = 565
$ws.Range("F45").Value = 78
$ws.Range("F46").Value = 490
$ws.Range("F47").Value = 539
$ws.Range("F48").Value = 251
$ws.Range("F49").Value = 171

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 163
$ws.Range("F16").Value = 145
$ws.Range("F21").Value = 4

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 834

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 834
$ws.Range("F3").Value = 760
$ws.Range("F5").Value = 2049
$ws.Range("F6").Value = 3521
$ws.Range("F7").Value = 61
$ws.Range("F8").Value = 1464
$ws.Range("F9").Value = 4876
$ws.Range("F10").Value = 1818
$ws.Range("F13").Value = 71
$ws.Range("F16").Value = 79
$ws.Range("F19").Value = 349
$ws.Range("F22").Value = 12
$ws.Range("F23").Value = 229
$ws.Range("F25").Value = 118
$ws.Range("F26").Value = 255
$ws.Range("F28").Value = 1009
$ws.Range("F29").Value = 38
$ws.Range("F30").Value = 1875
$ws.Range("F31").Value = 2335
$ws.Range("F33").Value = 1107
$ws.Range("F37").Value = 47
$ws.Range("F38").Value = 309
$ws.Range("F39").Value = 176
$ws.Range("F40").Value = 4
$ws.Range("F42").Value = 688
$ws.Range("F43").Value = 565
$ws.Range("F45").Value = 539
$ws.Range("F46").Value = 251
$ws.Range("F48").Value = 171
